$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value = 'ECs'
$ws.Range('B2').Value = 'Guca2a'
$ws.Range('C2').Value = 'Gucy2c'
$ws.Range('D2').Value = 'Inflammatory-Mac'
$ws.Range('E2').Value = 3
$ws.Range('F2').Value = 1
$ws.Range('G2').Value = 0.6424206666666666
$ws.Range('H2').Value = 1.927262
$ws.Range('I2').Value = 0.287181738380071
$ws.Range('J2').Value = 0.2977136782305576
$ws.Range('K2').Value = 3
$ws.Range('L2').Value = 1
$ws.Range('M2').Value = 0.574202
$ws.Range('N2').Value = 1.722606
$ws.Range('O2').Value = 0.944448982881442
$ws.Range('P2').Value = 0.9444489828814421
$ws.Range('Q2').Value = 0.3688792316413333
$ws.Range('R2').Value = 3.319913084772
$ws.Range('S2').Value = 0.2712285007151825
$ws.Range('T2').Value = 0.281175380594743

# Row 3
$ws.Range('A3').Value = 'ECs'
$ws.Range('B3').Value = 'Guca2a'
$ws.Range('C3').Value = 'Gucy2c'
$ws.Range('D3').Value = 'Resolving-Mac'
$ws.Range('E3').Value = 3
$ws.Range('F3').Value = 1
$ws.Range('G3').Value = 0.6424206666666666
$ws.Range('H3').Value = 1.927262
$ws.Range('I3').Value = 0.287181738380071
$ws.Range('J3').Value = 0.2977136782305576
$ws.Range('K3').Value = 1
$ws.Range('L3').Value = 0.3333333333333333
$ws.Range('M3').Value = 0.03377366666666667
$ws.Range('N3').Value = 0.101321
$ws.Range('O3').Value = 0.05555101711855792
$ws.Range('P3').Value = 0.05555101711855793
$ws.Range('Q3').Value = 0.02169690145577778
$ws.Range('R3').Value = 0.195272113102
$ws.Range('S3').Value = 0.01595323766488855
$ws.Range('T3').Value = 0.01653829763581455

# Row 4
$ws.Range('A4').Value = 'FAPs'
$ws.Range('B4').Value = 'Guca2a'
$ws.Range('C4').Value = 'Gucy2c'
$ws.Range('D4').Value = 'Inflammatory-Mac'
$ws.Range('E4').Value = 1
$ws.Range('F4').Value = 0.3333333333333333
$ws.Range('G4').Value = 0.3069103333333333
$ws.Range('H4').Value = 0.920731
$ws.Range('I4').Value = 0.1371983306682855
$ws.Range('J4').Value = 0.1422298642690508
$ws.Range('K4').Value = 3
$ws.Range('L4').Value = 1
$ws.Range('M4').Value = 0.574202
$ws.Range('N4').Value = 1.722606
$ws.Range('O4').Value = 0.944448982881442
$ws.Range('P4').Value = 0.9444489828814421
$ws.Range('Q4').Value = 0.1762285272206667
$ws.Range('R4').Value = 1.586056744986
$ws.Range('S4').Value = 0.129576823852694
$ws.Range('T4').Value = 0.1343288506442706

# Row 5
$ws.Range('A5').Value = 'FAPs'
$ws.Range('B5').Value = 'Guca2a'
$ws.Range('C5').Value = 'Gucy2c'
$ws.Range('D5').Value = 'Resolving-Mac'
$ws.Range('E5').Value = 1
$ws.Range('F5').Value = 0.3333333333333333
$ws.Range('G5').Value = 0.3069103333333333
$ws.Range('H5').Value = 0.920731
$ws.Range('I5').Value = 0.1371983306682855
$ws.Range('J5').Value = 0.1422298642690508
$ws.Range('K5').Value = 1
$ws.Range('L5').Value = 0.3333333333333333
$ws.Range('M5').Value = 0.03377366666666667
$ws.Range('N5').Value = 0.101321
$ws.Range('O5').Value = 0.05555101711855792
$ws.Range('P5').Value = 0.05555101711855793
$ws.Range('Q5').Value = 0.01036548729455556
$ws.Range('R5').Value = 0.09328938565099999
$ws.Range('S5').Value = 0.007621506815591497
$ws.Range('T5').Value = 0.007901013624780214

# Row 6
$ws.Range('A6').Value = 'Inflammatory-Mac'
$ws.Range('B6').Value = 'Guca2a'
$ws.Range('C6').Value = 'Gucy2c'
$ws.Range('D6').Value = 'Inflammatory-Mac'
$ws.Range('E6').Value = 2
$ws.Range('F6').Value = 0.6666666666666666
$ws.Range('G6').Value = 0.8184963333333334
$ws.Range('H6').Value = 2.455489
$ws.Range('I6').Value = 0.3658929608912241
$ws.Range('J6').Value = 0.3793115113797053
$ws.Range('K6').Value = 3
$ws.Range('L6').Value = 1
$ws.Range('M6').Value = 0.574202
$ws.Range('N6').Value = 1.722606
$ws.Range('O6').Value = 0.944448982881442
$ws.Range('P6').Value = 0.9444489828814421
$ws.Range('Q6').Value = 0.4699822315926667
$ws.Range('R6').Value = 4.229840084334
$ws.Range('S6').Value = 0.3455672347571959
$ws.Range('T6').Value = 0.3582403711177852

# Row 7
$ws.Range('A7').Value = 'Inflammatory-Mac'
$ws.Range('B7').Value = 'Guca2a'
$ws.Range('C7').Value = 'Gucy2c'
$ws.Range('D7').Value = 'Resolving-Mac'
$ws.Range('E7').Value = 2
$ws.Range('F7').Value = 0.6666666666666666
$ws.Range('G7').Value = 0.8184963333333334
$ws.Range('H7').Value = 2.455489
$ws.Range('I7').Value = 0.3658929608912241
$ws.Range('J7').Value = 0.3793115113797053
$ws.Range('K7').Value = 1
$ws.Range('L7').Value = 0.3333333333333333
$ws.Range('M7').Value = 0.03377366666666667
$ws.Range('N7').Value = 0.101321
$ws.Range('O7').Value = 0.05555101711855792
$ws.Range('P7').Value = 0.05555101711855793
$ws.Range('Q7').Value = 0.02764362232988889
$ws.Range('R7').Value = 0.248792600969
$ws.Range('S7').Value = 0.02032572613402824
$ws.Range('T7').Value = 0.02107114026192009

# Row 8
$ws.Range('A8').Value = 'MuSCs'
$ws.Range('B8').Value = 'Guca2a'
$ws.Range('C8').Value = 'Gucy2c'
$ws.Range('D8').Value = 'Inflammatory-Mac'
$ws.Range('E8').Value = 1
$ws.Range('F8').Value = 0.5
$ws.Range('G8').Value = 0.237407
$ws.Range('H8').Value = 0.474814
$ws.Range('I8').Value = 0.106128209289029
$ws.Range('J8').Value = 0.0733468632782486
$ws.Range('K8').Value = 3
$ws.Range('L8').Value = 1
$ws.Range('M8').Value = 0.574202
$ws.Range('N8').Value = 1.722606
$ws.Range('O8').Value = 0.944448982881442
$ws.Range('P8').Value = 0.9444489828814421
$ws.Range('Q8').Value = 0.136319574214
$ws.Range('R8').Value = 0.817917445284
$ws.Range('S8').Value = 0.1002326793180522
$ws.Range('T8').Value = 0.06927237042068608

# Row 9
$ws.Range('A9').Value = 'MuSCs'
$ws.Range('B9').Value = 'Guca2a'
$ws.Range('C9').Value = 'Gucy2c'
$ws.Range('D9').Value = 'Resolving-Mac'
$ws.Range('E9').Value = 1
$ws.Range('F9').Value = 0.5
$ws.Range('G9').Value = 0.237407
$ws.Range('H9').Value = 0.474814
$ws.Range('I9').Value = 0.106128209289029
$ws.Range('J9').Value = 0.0733468632782486
$ws.Range('K9').Value = 1
$ws.Range('L9').Value = 0.3333333333333333
$ws.Range('M9').Value = 0.03377366666666667
$ws.Range('N9').Value = 0.101321
$ws.Range('O9').Value = 0.05555101711855792
$ws.Range('P9').Value = 0.05555101711855793
$ws.Range('Q9').Value = 0.008018104882333333
$ws.Range('R9').Value = 0.048108629294
$ws.Range('S9').Value = 0.005895529970976748
$ws.Range('T9').Value = 0.004074492857562516

# Row 10
$ws.Range('A10').Value = 'Resolving-Mac'
$ws.Range('B10').Value = 'Guca2a'
$ws.Range('C10').Value = 'Gucy2c'
$ws.Range('D10').Value = 'Inflammatory-Mac'
$ws.Range('E10').Value = 1
$ws.Range('F10').Value = 0.3333333333333333
$ws.Range('G10').Value = 0.2317486666666667
$ws.Range('H10').Value = 0.695246
$ws.Range('I10').Value = 0.1035987607713901
$ws.Range('J10').Value = 0.1073980828424377
$ws.Range('K10').Value = 3
$ws.Range('L10').Value = 1
$ws.Range('M10').Value = 0.574202
$ws.Range('N10').Value = 1.722606
$ws.Range('O10').Value = 0.944448982881442
$ws.Range('P10').Value = 0.9444489828814421
$ws.Range('Q10').Value = 0.1330705478973333
$ws.Range('R10').Value = 1.197634931076
$ws.Range('S10').Value = 0.09784374423831725
$ws.Range('T10').Value = 0.1014320101039572

# Row 11
$ws.Range('A11').Value = 'Resolving-Mac'
$ws.Range('B11').Value = 'Guca2a'
$ws.Range('C11').Value = 'Gucy2c'
$ws.Range('D11').Value = 'Resolving-Mac'
$ws.Range('E11').Value = 1
$ws.Range('F11').Value = 0.3333333333333333
$ws.Range('G11').Value = 0.2317486666666667
$ws.Range('H11').Value = 0.695246
$ws.Range('I11').Value = 0.1035987607713901
$ws.Range('J11').Value = 0.1073980828424377
$ws.Range('K11').Value = 1
$ws.Range('L11').Value = 0.3333333333333333
$ws.Range('M11').Value = 0.03377366666666667
$ws.Range('N11').Value = 0.101321
$ws.Range('O11').Value = 0.05555101711855792
$ws.Range('P11').Value = 0.05555101711855793
$ws.Range('Q11').Value = 0.007827002218444444
$ws.Range('R11').Value = 0.070443019966
$ws.Range('S11').Value = 0.005755016533072881
$ws.Range('T11').Value = 0.005966072738480561
